$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = [string]$cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $newParts = $parts[1..($parts.Count - 1)] + $parts[0]
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
